# NYPD CompStat weekly report update: new crime data collected.
# Updates the "Volume/Number" and "Report Covering the Week" header text,
# and refreshes the crime-statistics grid (rows 15-30) with the latest
# Week-to-Date / 28-Day / Year-to-Date / 2-Year figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text (report issue number + covered week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  35"
$ws.Range("C9").Value = "Report Covering the Week  8/28/2023  Through  9/3/2023"

# Stable donor cells already carrying the placeholder text styling used
# throughout the grid for "no data" (shared text "0") and "not applicable"
# (shared text "***.*") values - row 14 (Murder) is untouched by this
# week's edits, so it is safe to copy format+value from here.
$txtZero = $ws.Range("C14")
$txtNA = $ws.Range("E14")

# --- Numeric cell updates ---
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = -48
$ws.Range("I16").Value = 124
$ws.Range("J16").Value = 165
$ws.Range("K16").Value = -24.848484848484
$ws.Range("L16").Value = 11.711711711711
$ws.Range("M16").Value = 39.325842696629
$ws.Range("N16").Value = -78.397212543554
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -83.333333333333
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = -56.666666666666
$ws.Range("I17").Value = 146
$ws.Range("J17").Value = 131
$ws.Range("K17").Value = 11.450381679389
$ws.Range("L17").Value = 22.689075630252
$ws.Range("M17").Value = 135.483870967742
$ws.Range("N17").Value = -29.807692307692
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 12
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 40
$ws.Range("H18").Value = -70
$ws.Range("I18").Value = 194
$ws.Range("J18").Value = 306
$ws.Range("K18").Value = -36.601307189542
$ws.Range("L18").Value = 45.864661654135
$ws.Range("M18").Value = 52.755905511811
$ws.Range("N18").Value = -65.602836879432
$ws.Range("C19").Value = 27
$ws.Range("D19").Value = 24
$ws.Range("E19").Value = 12.5
$ws.Range("G19").Value = 102
$ws.Range("H19").Value = 7.843137254901
$ws.Range("I19").Value = 847
$ws.Range("J19").Value = 880
$ws.Range("K19").Value = -3.75
$ws.Range("L19").Value = 92.5
$ws.Range("M19").Value = 17.638888888888
$ws.Range("N19").Value = -48.635536688902
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 14.285714285714
$ws.Range("J20").Value = 44
$ws.Range("K20").Value = -27.272727272727
$ws.Range("L20").Value = -3.030303030303
$ws.Range("N20").Value = -93.456032719836
$ws.Range("C21").Value = 39
$ws.Range("D21").Value = 53
$ws.Range("E21").Value = -26.415094339622
$ws.Range("F21").Value = 157
$ws.Range("G21").Value = 204
$ws.Range("H21").Value = -23.039215686274
$ws.Range("I21").Value = 1351
$ws.Range("J21").Value = 1538
$ws.Range("K21").Value = -12.158647594278
$ws.Range("L21").Value = 60.260972716488
$ws.Range("M21").Value = 31.037827352085
$ws.Range("N21").Value = -61.322645290581
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = -33.333333333333
$ws.Range("J22").Value = 29
$ws.Range("K22").Value = -3.448275862068
$ws.Range("C24").Value = 44
$ws.Range("D24").Value = 42
$ws.Range("E24").Value = 4.761904761904
$ws.Range("F24").Value = 181
$ws.Range("G24").Value = 169
$ws.Range("H24").Value = 7.100591715976
$ws.Range("I24").Value = 1437
$ws.Range("J24").Value = 1332
$ws.Range("K24").Value = 7.882882882882
$ws.Range("L24").Value = 71.275327771156
$ws.Range("M24").Value = 39.514563106796
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 125
$ws.Range("F25").Value = 27
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = -22.857142857142
$ws.Range("I25").Value = 289
$ws.Range("J25").Value = 270
$ws.Range("K25").Value = 7.037037037037
$ws.Range("L25").Value = 48.205128205128
$ws.Range("M25").Value = 71.005917159763
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 45
$ws.Range("J27").Value = 47
$ws.Range("K27").Value = -4.255319148936
$ws.Range("L27").Value = 15.384615384615
$ws.Range("F30").Value = 2
$ws.Range("L30").Value = -45.454545454545

# --- Text placeholder cell updates (copy format+value from stable donor cells) ---
$txtZero.Copy($ws.Range("C15"))
$txtZero.Copy($ws.Range("G15"))
$txtNA.Copy($ws.Range("H15"))
$txtZero.Copy($ws.Range("C20"))
$txtZero.Copy($ws.Range("C26"))
$txtZero.Copy($ws.Range("G26"))
$txtNA.Copy($ws.Range("H26"))
$txtZero.Copy($ws.Range("C30"))
$txtZero.Copy($ws.Range("G30"))
$txtNA.Copy($ws.Range("H30"))
